# LWP2_0015_lab_timing.xlsx - "update segmentation with new timing"
#
# The lab re-segmented the in-lab session: the R1/T1/E1/S1/... event
# timestamps in column B (rows 7-24 of the "In Lab" sheet) are replaced
# with freshly re-timed values, and reformatted to show seconds
# (h:mm:ss) like the rest of the sheet instead of the old h:mm format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("In Lab")

# New timing values for the event markers (row -> new Excel time serial)
$newTimes = @{
    7  = 0.60646990740740747   # R1
    8  = 0.61371527777777779   # T1
    9  = 0.61527777777777781   # E1  (value unchanged, only re-formatted)
    10 = 0.61636574074074069   # S1
    11 = 0.62479166666666663   # T2
    12 = 0.62649305555555557   # E2
    13 = 0.6274305555555556    # R2
    14 = 0.63484953703703706   # T3
    15 = 0.63658564814814811   # E3
    16 = 0.63739583333333327   # S2
    17 = 0.64008101851851851   # R3
    18 = 0.64753472222222219   # T4
    19 = 0.64887731481481481   # E4
    20 = 0.64943287037037034   # S3
    21 = 0.6521527777777778    # S4
    22 = 0.65413194444444445   # T5
    23 = 0.65567129629629628   # E5
    24 = 0.65649305555555559   # S5
}

foreach ($row in $newTimes.Keys) {
    $cell = $ws.Cells.Item($row, 2)   # column B
    $cell.Value = $newTimes[$row]
    $cell.NumberFormat = "h:mm:ss"
}

# Row 11 (T2) also picks up right-aligned text, matching the formatting
# used elsewhere on timestamps that need to line up.
$ws.Range("B11").HorizontalAlignment = -4152   # xlRight

# The view was scrolled/re-selected while reviewing the new segmentation.
$ws.Application.Goto($ws.Range("A5"), $true)
$ws.Range("B25").Select()
